$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 69; this shifts rows 69..164 down to 70..165,
# preserving all existing formatting (e.g. the date style on column D).
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new record.
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value = "Los Lagos"
$ws.Range("D69").Value = 44495
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 100112043
$ws.Range("G69").Value = "Pepino ensalada"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 400
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = 12000
$ws.Range("N69").Value = "$/caja 60 unidades"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 200
$ws.Range("Q69").Value = 60
$ws.Range("R69").Value = "Hortaliza"
